$d = $word.ActiveDocument

$pairs = @(
    @("49×19=", "56×32="),
    @("54×91=", "16×98="),
    @("83×90=", "65×22="),
    @("72×31=", "31×79="),
    @("11×60=", "85×97="),
    @("88×62=", "55×47="),
    @("71×75=", "84×43="),
    @("85×84=", "45×11="),
    @("65×36=", "88×34="),
    @("97×78=", "30×37="),
    @("51×13=", "37×68="),
    @("64×87=", "23×17="),
    @("55×71=", "15×50="),
    @("70×32=", "53×83="),
    @("52×70=", "37×81="),
    @("53×92=", "49×54="),
    @("50×99=", "28×37="),
    @("43×77=", "91×52="),
    @("20×64=", "35×52="),
    @("55×19=", "96×26="),
    @("50×46=", "58×57="),
    @("87×60=", "88×55="),
    @("39×36=", "39×87="),
    @("15×87=", "20×32="),
    @("72×48=", "66×21=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}
